$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("M15").Value = -1277.8361
$ws.Range("K15").Value = 1446.8361
$ws.Range("H15").Value = 482.2787
$ws.Range("I15").Value = 482.2787

# Row 19
$ws.Range("I19").Value = 1127.5
$ws.Range("H19").Value = 1571.9546
$ws.Range("M19").Value = -952.5
$ws.Range("K19").Value = 1127.5

# Row 87
$ws.Range("N87").ClearContents()
$ws.Range("H87").Value = 45000
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0

# Row 90
$ws.Range("H90").Value = 45000
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("J90").Value = 0

# Row 92
$ws.Range("I92").Value = 639.5
$ws.Range("M92").Value = 608.5
$ws.Range("H92").Value = 917.1818
$ws.Range("K92").Value = 639.5

# Row 100
$ws.Range("I100").Value = 1261.5
$ws.Range("M100").Value = -720.5
$ws.Range("K100").Value = 1261.5
$ws.Range("H100").Value = 1577.4445

# Row 111
$ws.Range("H111").Value = 691.65216
$ws.Range("M111").Value = 1182.47065
$ws.Range("K111").Value = 1884.52935
$ws.Range("I111").Value = 628.17645

# Row 137
$ws.Range("I137").Value = 2002.7222
$ws.Range("M137").Value = -3458.1666
$ws.Range("K137").Value = 6008.1666
$ws.Range("H137").Value = 405631.6

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("K32").Value = 2353.6538
$ws.Range("I32").Value = 2353.6538
$ws.Range("H32").Value = 2585.5334
$ws.Range("M32").Value = -2066.6538

# Row 97
$ws.Range("M97").Value = -476.5
$ws.Range("K97").Value = 972.5
$ws.Range("I97").Value = 972.5
$ws.Range("H97").Value = 1051.5385

# Row 110
$ws.Range("I110").Value = 1431
$ws.Range("K110").Value = 1431
$ws.Range("M110").Value = 614
$ws.Range("H110").Value = 2012.8334

# Row 122
$ws.Range("H122").Value = 1552.2
$ws.Range("J122").Value = 2232
$ws.Range("N122").Value = -11596
$ws.Range("L122").Value = 6696
$ws.Range("K122").Value = 4146.75
$ws.Range("I122").Value = 1382.25
$ws.Range("M122").Value = -1696.75

# Row 132
$ws.Range("L132").Value = 9498.999899999999
$ws.Range("I132").Value = 3231
$ws.Range("H132").Value = 3220.7896
$ws.Range("N132").Value = -14558.9999
$ws.Range("J132").Value = 3166.3333
$ws.Range("K132").Value = 9693
$ws.Range("M132").Value = -7163

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("I94").Value = 1749.7693
$ws.Range("M94").Value = -1298.7693
$ws.Range("K94").Value = 1749.7693
$ws.Range("L94").Value = 2457.5
$ws.Range("N94").Value = -3359.5
$ws.Range("J94").Value = 2457.5
$ws.Range("H94").Value = 1973.2632

# Row 107
$ws.Range("H107").Value = 2111.111
$ws.Range("K107").Value = 2021.409
$ws.Range("M107").Value = -101.4090000000001
$ws.Range("I107").Value = 2021.409

# Row 134
$ws.Range("K134").Value = 8578.7775
$ws.Range("H134").Value = 4519.6665
$ws.Range("L134").Value = 28499.667
$ws.Range("I134").Value = 2859.5925
$ws.Range("M134").Value = -6043.7775
$ws.Range("N134").Value = -33569.667
$ws.Range("J134").Value = 9499.888999999999

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("I7").Value = 51.25
$ws.Range("K7").Value = 51.25
$ws.Range("L7").Value = 224.85715
$ws.Range("M7").Value = 61.75
$ws.Range("J7").Value = 224.85715
$ws.Range("N7").Value = -450.85715
$ws.Range("H7").Value = 176.96552

# Row 50
$ws.Range("N50").ClearContents()
$ws.Range("L50").Value = 0
$ws.Range("H50").Value = 0
$ws.Range("J50").Value = 0

# Row 51
$ws.Range("J51").Value = 23222.111
$ws.Range("L51").Value = 23222.111
$ws.Range("N51").Value = -24694.111
$ws.Range("H51").Value = 31213.066

# Row 61
$ws.Range("L61").Value = 23222.111
$ws.Range("J61").Value = 23222.111
$ws.Range("H61").Value = 31213.066
$ws.Range("N61").Value = -23918.111

# Row 132
$ws.Range("H132").Value = 1585327.6

$ws = $wb.Worksheets.Item("CUL")
# Row 6
$ws.Range("H6").Value = 166967.17
$ws.Range("K6").Value = 500901.51
$ws.Range("I6").Value = 166967.17
$ws.Range("M6").Value = -500788.51

# Row 131
$ws.Range("I131").Value = 759.2
$ws.Range("L131").Value = 6171.6666
$ws.Range("N131").Value = -16251.6666
$ws.Range("M131").Value = 2762.4
$ws.Range("J131").Value = 2057.2222
$ws.Range("K131").Value = 2277.6
$ws.Range("H131").Value = 1593.6428

# Row 133
$ws.Range("J133").Value = 9999
$ws.Range("H133").Value = 5674.857
$ws.Range("L133").Value = 29997
$ws.Range("N133").Value = -40117

# Row 137
$ws.Range("L137").Value = 28599.501
$ws.Range("I137").Value = 2230.111
$ws.Range("J137").Value = 9533.166999999999
$ws.Range("M137").Value = -1590.333
$ws.Range("N137").Value = -38799.501
$ws.Range("K137").Value = 6690.333
$ws.Range("H137").Value = 4055.875

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("I2").Value = 109.42857
$ws.Range("K2").Value = 109.42857
$ws.Range("H2").Value = 168.42105
$ws.Range("J2").Value = 333.6
$ws.Range("M2").Value = 3.571430000000007
$ws.Range("N2").Value = -559.6
$ws.Range("L2").Value = 333.6

# Row 97
$ws.Range("M97").Value = -2434.8572
$ws.Range("K97").Value = 2930.8572
$ws.Range("I97").Value = 2930.8572
$ws.Range("H97").Value = 2973.8572

# Row 122
$ws.Range("H122").Value = 77863.664
$ws.Range("J122").Value = 5622.125
$ws.Range("N122").Value = -21766.375
$ws.Range("L122").Value = 16866.375

# Row 134
$ws.Range("H134").Value = 28903.285
$ws.Range("L134").Value = 86709.855
$ws.Range("N134").Value = -91779.855
$ws.Range("J134").Value = 28903.285

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("I7").Value = 15416.5
$ws.Range("K7").Value = 15416.5
$ws.Range("L7").Value = 4486.875
$ws.Range("M7").Value = -15304.5
$ws.Range("J7").Value = 4486.875
$ws.Range("N7").Value = -4710.875
$ws.Range("H7").Value = 9951.6875

# Row 22
$ws.Range("M22").Value = -2355.6843
$ws.Range("I22").Value = 2650.6843
$ws.Range("K22").Value = 2650.6843
$ws.Range("H22").Value = 73863.07000000001

# Row 27
$ws.Range("M27").Value = -2543.6843
$ws.Range("I27").Value = 2650.6843
$ws.Range("H27").Value = 73863.07000000001
$ws.Range("K27").Value = 2650.6843

# Row 55
$ws.Range("J55").Value = 3399.4
$ws.Range("I55").Value = 1711.3077
$ws.Range("M55").Value = -1538.3077
$ws.Range("L55").Value = 3399.4
$ws.Range("N55").Value = -3745.4
$ws.Range("K55").Value = 1711.3077
$ws.Range("H55").Value = 2445.261

# Row 126
$ws.Range("I126").Value = 15416.5
$ws.Range("N126").Value = -18400.625
$ws.Range("J126").Value = 4486.875
$ws.Range("L126").Value = 13460.625
$ws.Range("K126").Value = 46249.5
$ws.Range("M126").Value = -43779.5
$ws.Range("H126").Value = 9951.6875

# Row 132
$ws.Range("I132").Value = 1000
$ws.Range("H132").Value = 6334.25
$ws.Range("K132").Value = 3000
$ws.Range("M132").Value = -470

# Row 133
$ws.Range("J133").Value = 39750
$ws.Range("H133").Value = 39750
$ws.Range("L133").Value = 39750
$ws.Range("N133").Value = -44810

# Row 136
$ws.Range("N136").Value = -14858.667
$ws.Range("K136").Value = 15511.092
$ws.Range("M136").Value = -12961.092
$ws.Range("J136").Value = 3252.889
$ws.Range("L136").Value = 9758.667000000001
$ws.Range("I136").Value = 5170.364
$ws.Range("H136").Value = 4307.5

$ws = $wb.Worksheets.Item("WVR")
# Row 27
$ws.Range("L27").Value = 64750
$ws.Range("N27").Value = -64888
$ws.Range("J27").Value = 64750
$ws.Range("H27").Value = 64750

# Row 31
$ws.Range("H31").Value = 0
$ws.Range("K31").Value = 0
$ws.Range("I31").Value = 0
$ws.Range("M31").ClearContents()

# Row 108
$ws.Range("J108").Value = 100000
$ws.Range("L108").Value = 100000
$ws.Range("N108").Value = -107680
$ws.Range("H108").Value = 100000

# Row 109
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52774
$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000

# Row 110
$ws.Range("L110").Value = 20000
$ws.Range("J110").Value = 20000
$ws.Range("N110").Value = -28180
$ws.Range("H110").Value = 20000

# Row 111
$ws.Range("H111").Value = 25551.25
$ws.Range("N111").Value = -28697.5
$ws.Range("L111").Value = 20517.5
$ws.Range("J111").Value = 20517.5

# Row 114
$ws.Range("J114").Value = 67500
$ws.Range("N114").Value = -76178
$ws.Range("L114").Value = 67500
$ws.Range("H114").Value = 67500

# Row 115
$ws.Range("L115").Value = 0
$ws.Range("H115").Value = 0
$ws.Range("J115").Value = 0
$ws.Range("N115").ClearContents()
